$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44329
$ws.Cells.Item(2, 10).Value = 40

# Row 3
$ws.Cells.Item(3, 4).Value = 44455
$ws.Cells.Item(3, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(3, 15).Value = 'Perú'
$ws.Cells.Item(3, 16).Value = 1000
$ws.Cells.Item(3, 17).Value = 20

# Row 4
$ws.Cells.Item(4, 4).Value = 44441
$ws.Cells.Item(4, 10).Value = 40
$ws.Cells.Item(4, 15).Value = 'Perú'

# Row 5
$ws.Cells.Item(5, 4).Value = 44179

# Row 6
$ws.Cells.Item(6, 4).Value = 44341
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(6, 11).Value = 17000
$ws.Cells.Item(6, 12).Value = 18000
$ws.Cells.Item(6, 13).Value = 17500
$ws.Cells.Item(6, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(6, 15).Value = 'Perú'
$ws.Cells.Item(6, 16).Value = 875
$ws.Cells.Item(6, 17).Value = 20

# Row 7
$ws.Cells.Item(7, 4).Value = 44340
$ws.Cells.Item(7, 10).Value = 40
$ws.Cells.Item(7, 11).Value = 18000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 18000
$ws.Cells.Item(7, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(7, 15).Value = 'Perú'
$ws.Cells.Item(7, 16).Value = 900
$ws.Cells.Item(7, 17).Value = 20

# Row 8
$ws.Cells.Item(8, 4).Value = 44175
$ws.Cells.Item(8, 10).Value = 20
$ws.Cells.Item(8, 15).Value = 'Región de Arica y Parinacota'

# Row 11
$ws.Cells.Item(11, 4).Value = 44424
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(11, 16).Value = 1333
$ws.Cells.Item(11, 17).Value = 15

# Row 12
$ws.Cells.Item(12, 4).Value = 44315
$ws.Cells.Item(12, 10).Value = 30
$ws.Cells.Item(12, 11).Value = 20000
$ws.Cells.Item(12, 12).Value = 20000
$ws.Cells.Item(12, 13).Value = 20000
$ws.Cells.Item(12, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(12, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(12, 16).Value = 1333
$ws.Cells.Item(12, 17).Value = 15

# Row 13
$ws.Cells.Item(13, 4).Value = 44315
$ws.Cells.Item(13, 10).Value = 30
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 20000
$ws.Cells.Item(13, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(13, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(13, 16).Value = 1000
$ws.Cells.Item(13, 17).Value = 20

# Row 14
$ws.Cells.Item(14, 4).Value = 44294
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(14, 16).Value = 1333
$ws.Cells.Item(14, 17).Value = 15

# Row 15
$ws.Cells.Item(15, 4).Value = 44316
$ws.Cells.Item(15, 10).Value = 20
$ws.Cells.Item(15, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(15, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(15, 16).Value = 1333
$ws.Cells.Item(15, 17).Value = 15

# Row 17
$ws.Cells.Item(17, 4).Value = 44186
$ws.Cells.Item(17, 10).Value = 20
$ws.Cells.Item(17, 11).Value = 20000
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = 20000
$ws.Cells.Item(17, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(17, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(17, 16).Value = 1333
$ws.Cells.Item(17, 17).Value = 15

# Row 19
$ws.Cells.Item(19, 4).Value = 44385
$ws.Cells.Item(19, 10).Value = 18
$ws.Cells.Item(19, 15).Value = 'Región de Arica y Parinacota'

# Row 20
$ws.Cells.Item(20, 4).Value = 44452
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(20, 15).Value = 'Perú'
$ws.Cells.Item(20, 16).Value = 1000
$ws.Cells.Item(20, 17).Value = 20

# Row 21
$ws.Cells.Item(21, 4).Value = 44369
$ws.Cells.Item(21, 10).Value = 20

# Row 22
$ws.Cells.Item(22, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(22, 16).Value = 1000
$ws.Cells.Item(22, 17).Value = 20

# Row 23
$ws.Cells.Item(23, 4).Value = 44321
$ws.Cells.Item(23, 10).Value = 15
$ws.Cells.Item(23, 11).Value = 25000
$ws.Cells.Item(23, 12).Value = 25000
$ws.Cells.Item(23, 13).Value = 25000
$ws.Cells.Item(23, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(23, 15).Value = 'Perú'
$ws.Cells.Item(23, 16).Value = 1667
$ws.Cells.Item(23, 17).Value = 15

# Row 24
$ws.Cells.Item(24, 4).Value = 44425
$ws.Cells.Item(24, 10).Value = 10

# Row 25
$ws.Cells.Item(25, 4).Value = 44466
$ws.Cells.Item(25, 10).Value = 20
$ws.Cells.Item(25, 11).Value = 25000
$ws.Cells.Item(25, 12).Value = 25000
$ws.Cells.Item(25, 13).Value = 25000
$ws.Cells.Item(25, 15).Value = 'Perú'
$ws.Cells.Item(25, 16).Value = 1667

# Row 26
$ws.Cells.Item(26, 4).Value = 44438
$ws.Cells.Item(26, 10).Value = 40

# Row 27
$ws.Cells.Item(27, 4).Value = 44389
$ws.Cells.Item(27, 10).Value = 45
